$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.178289532661438
$ws.Range("B1").Value = 2.375503063201904
$ws.Range("C1").Value = 3.529220819473267
$ws.Range("D1").Value = 1.834197759628296
$ws.Range("E1").Value = 1.208972692489624
